$wb = $excel.ActiveWorkbook

# --- Sheet "Trees": add a new row (row 6) ---
$trees = $wb.Worksheets.Item("Trees")

$trees.Cells.Item(6,1).Value = "tang"
$trees.Cells.Item(6,2).Value = "testingtree"
$trees.Cells.Item(6,3).Value = 111
$trees.Cells.Item(6,4).Value = 11

# E6 is a tapping date - copy the date format from an existing date cell (E2)
# then set the value, so the style matches the existing date style exactly.
$trees.Cells.Item(2,5).Copy()
$trees.Cells.Item(6,5).PasteSpecial(-4122)
$trees.Cells.Item(6,5).Value = 44318

$trees.Cells.Item(6,6).Value = "N/A"
$trees.Cells.Item(6,7).Value = "N/A"
$trees.Cells.Item(6,10).Value = "testing"
$trees.Cells.Item(6,11).Value = "N/A"

# --- Sheet "Seasons": add a new row (row 7) ---
$seasons = $wb.Worksheets.Item("Seasons")

$seasons.Cells.Item(7,1).Value = "tang"
$seasons.Cells.Item(7,2).Value = "testingtree"
$seasons.Cells.Item(7,3).Value = 2021

# D7 is a date - copy the date format from an existing date cell (D2)
$seasons.Cells.Item(2,4).Copy()
$seasons.Cells.Item(7,4).PasteSpecial(-4122)
$seasons.Cells.Item(7,4).Value = 44318

$seasons.Cells.Item(7,5).Value = "N/A"
$seasons.Cells.Item(7,6).Value = "N/A"
$seasons.Cells.Item(7,7).Value = "testing"
$seasons.Cells.Item(7,8).Value = "N/A"

# --- Sheet "Saps": delete row 12 (SapID 607e0e31122a9d2f1092b967) ---
$saps = $wb.Worksheets.Item("Saps")
$saps.Rows.Item(12).Delete()
